$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the <w:proofErr spellStart/spellEnd> wrapper around the
#    "ИНБс" run (first occurrence, in the title block) while keeping
#    the run text/formatting (rStyle=spellingerror) intact.
#    Word's Find/Replace engine only drops proofErr markers that fall
#    strictly inside a replaced span, so we replace a span that spans
#    from just before "ИНБс" to just after it (merging the 3 runs into
#    one run using the leading run's formatting), and then re-split
#    the "ИНБс" substring back out into its own run by re-applying its
#    original character style.
# ------------------------------------------------------------------
$rngMerge = $d.Content
$rngMerge.Find.Execute(" ИНБс ", $true, $false, $false, $false, $false, $true, 1, $false, " ИНБс ", 2)

$rngRestyle = $d.Content
$rngRestyle.Find.Execute("ИНБс", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngRestyle.Style = "spellingerror"

# ------------------------------------------------------------------
# 2) Mark the three figure-holding runs (the ones whose paragraph
#    contains a picture/drawing) as NoProofing -> emits <w:noProof/>.
#    The very first inline picture already has it in the source doc,
#    so only shapes 2..N need the flag flipped on.
# ------------------------------------------------------------------
for ($i = 2; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# ------------------------------------------------------------------
# 3) Delete the "empty paragraph" + "Расчет вручную:" heading +
#    the long "Начиная с начала массива..." explanation paragraph
#    that used to sit between "Верификация результатов:" and
#    "Дан массив из 10 чисел...".
# ------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(62)
$pEnd = $d.Paragraphs.Item(64)
$rngCut = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rngCut.Delete()

Write-Output "done"
